$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 14, pushing existing rows 14-22 down to 15-23.
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new weekly data entry.
$ws.Range("A14").Value = 8
$ws.Range("B14").Value = "Terminal La Palmera de La Serena"
$ws.Range("C14").Value = "Coquimbo"
$ws.Range("D14").Value = 44634
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100101
$ws.Range("H14").Value = "Berries"
$ws.Range("I14").Value = 100101001
$ws.Range("J14").Value = "Arándano (blue)"
$ws.Range("K14").Value = "Sin especificar"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 160
$ws.Range("N14").Value = 6000
$ws.Range("O14").Value = 6500
$ws.Range("P14").Value = 6250
$ws.Range("Q14").Value = "`$/bandeja 2 kilos"
$ws.Range("R14").Value = "Provincia de Linares"
$ws.Range("S14").Value = 3125
$ws.Range("T14").Value = 2
